# Validation: corrected error in runall.bat. Updated values in V&V comparisons
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("J3").Value = 1.1399999999999999
$ws.Range("K3").Formula = "=0.48/2"

# Row 5
$ws.Range("J5").Value = 1.23

# Row 8
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Formula = "=0.61/2"

# Row 9
$ws.Range("J9").Value = 0.85
$ws.Range("K9").Formula = "=0.49/2"

# Row 15
$ws.Range("J15").Value = 1.0900000000000001
$ws.Range("K15").Formula = "=0.93/2"

# Update the active cell / selection to match the author's final cursor position
$ws.Range("K16").Select()
